# Insert a "property" / "value" header row at the top of every sheet
# (carbs, potatoes, pasta), pushing the existing name/value rows down by
# one, and leave the "pasta" sheet selected with the new header row
# highlighted.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Rows.Item(1).Insert()
    $ws.Cells.Item(1, 1).Value = "property"
    $ws.Cells.Item(1, 2).Value = "value"
    $ws.Range("A1:B1").Select()
}

$wb.Worksheets.Item("pasta").Activate()
